$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column S with data for 2022, mirroring the formatting of column R
# Row 4 header (year 2022)
$ws.Range("R4").Copy($ws.Range("S4")) | Out-Null
$ws.Range("S4").Value = 2022

# Row 5 (share of renewable energy)
$ws.Range("R5").Copy($ws.Range("S5")) | Out-Null
$ws.Range("S5").Value = 30

# Row 6 (electricity production)
$ws.Range("R6").Copy($ws.Range("S6")) | Out-Null
$ws.Range("S6").Value = 11928.6

# Update the active selection to T3 as in the target workbook
$ws.Range("T3").Select() | Out-Null
